{"js": "// Update the date line and every two-digit multiplication equation in the\n// table with the newly generated values.\nconst replacements = [\n  [\"2023-12-19 Tuesday\", \"2023-12-20 Wednesday\"],\n  [\"92\u00d724=2208\", \"47\u00d782=3854\"],\n  [\"92\u00d798=9016\", \"12\u00d770=840\"],\n  [\"86\u00d723=1978\", \"70\u00d751=3570\"],\n  [\"34\u00d755=1870\", \"33\u00d726=858\"],\n  [\"85\u00d792=7820\", \"43\u00d794=4042\"],\n  [\"99\u00d796=9504\", \"19\u00d730=570\"],\n  [\"22\u00d757=1254\", \"56\u00d767=3752\"],\n  [\"80\u00d746=3680\", \"83\u00d761=5063\"],\n  [\"34\u00d761=2074\", \"12\u00d720=240\"],\n  [\"28\u00d796=2688\", \"20\u00d770=1400\"],\n  [\"51\u00d737=1887\", \"67\u00d742=2814\"],\n  [\"87\u00d725=2175\", \"45\u00d777=3465\"],\n  [\"72\u00d759=4248\", \"25\u00d780=2000\"],\n  [\"19\u00d765=1235\", \"76\u00d734=2584\"],\n  [\"46\u00d736=1656\", \"75\u00d751=3825\"],\n  [\"64\u00d799=6336\", \"26\u00d793=2418\"],\n  [\"12\u00d787=1044\", \"45\u00d757=2565\"],\n  [\"65\u00d776=4940\", \"58\u00d713=754\"],\n  [\"34\u00d764=2176\", \"84\u00d718=1512\"],\n  [\"54\u00d760=3240\", \"63\u00d728=1764\"],\n  [\"29\u00d716=464\", \"50\u00d756=2800\"],\n  [\"42\u00d779=3318\", \"43\u00d734=1462\"],\n  [\"77\u00d720=1540\", \"71\u00d748=3408\"],\n  [\"80\u00d735=2800\", \"35\u00d776=2660\"],\n  [\"69\u00d789=6141\", \"72\u00d783=5976\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit multiplication equation in the\n# table with the newly generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-12-19 Tuesday\", \"2023-12-20 Wednesday\"),\n    @(\"92\u00d724=2208\", \"47\u00d782=3854\"),\n    @(\"92\u00d798=9016\", \"12\u00d770=840\"),\n    @(\"86\u00d723=1978\", \"70\u00d751=3570\"),\n    @(\"34\u00d755=1870\", \"33\u00d726=858\"),\n    @(\"85\u00d792=7820\", \"43\u00d794=4042\"),\n    @(\"99\u00d796=9504\", \"19\u00d730=570\"),\n    @(\"22\u00d757=1254\", \"56\u00d767=3752\"),\n    @(\"80\u00d746=3680\", \"83\u00d761=5063\"),\n    @(\"34\u00d761=2074\", \"12\u00d720=240\"),\n    @(\"28\u00d796=2688\", \"20\u00d770=1400\"),\n    @(\"51\u00d737=1887\", \"67\u00d742=2814\"),\n    @(\"87\u00d725=2175\", \"45\u00d777=3465\"),\n    @(\"72\u00d759=4248\", \"25\u00d780=2000\"),\n    @(\"19\u00d765=1235\", \"76\u00d734=2584\"),\n    @(\"46\u00d736=1656\", \"75\u00d751=3825\"),\n    @(\"64\u00d799=6336\", \"26\u00d793=2418\"),\n    @(\"12\u00d787=1044\", \"45\u00d757=2565\"),\n    @(\"65\u00d776=4940\", \"58\u00d713=754\"),\n    @(\"34\u00d764=2176\", \"84\u00d718=1512\"),\n    @(\"54\u00d760=3240\", \"63\u00d728=1764\"),\n    @(\"29\u00d716=464\", \"50\u00d756=2800\"),\n    @(\"42\u00d779=3318\", \"43\u00d734=1462\"),\n    @(\"77\u00d720=1540\", \"71\u00d748=3408\"),\n    @(\"80\u00d735=2800\", \"35\u00d776=2660\"),\n    @(\"69\u00d789=6141\", \"72\u00d783=5976\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
